{"js": "// The document repeats a \"2018 observation dates for the Perseus\n// constellation\" sentence in four places. Each occurrence is replaced with\n// a single, unformatted run containing the new (translated/updated) Orion\n// sentence \u2014 the \"2018 \" lead-in is dropped and the constellation / date\n// range text is swapped for the Orion one.\nconst OLD_TAIL = \"30 \u039f\u03ba\u03c4\u03c9\u03b2\u03c1\u03af\u03bf\u03c5-8 \u039d\u03bf\u03b5\u03bc\u03b2\u03c1\u03af\u03bf\u03c5 \u03ba\u03b1\u03b9 29 \u039d\u03bf\u03b5\u03bc\u03b2\u03c1\u03af\u03bf\u03c5-8 \u0394\u03b5\u03ba\u03b5\u03bc\u03b2\u03c1\u03af\u03bf\u03c5\";\nconst NEW_TEXT =\n  \"\u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2 \u03c0\u03b1\u03c1\u03b1\u03c4\u03ae\u03c1\u03b7\u03c3\u03b7\u03c2 \u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 Orion: 16-25 \u0399\u03b1\u03bd\u03bf\u03c5\u03b1\u03c1\u03af\u03bf\u03c5, 14-23 \u03a6\u03b5\u03b2\u03c1\u03bf\u03c5\u03b1\u03c1\u03af\u03bf\u03c5, 14-24 \u039c\u03b1\u03c1\u03c4\u03af\u03bf\u03c5\";\n\nconst body = context.document.body;\n\n// Locate every paragraph that still holds the old (Perseus / 2018) dates\n// sentence by searching for a distinctive tail substring shared by all of\n// them, then resolving each hit back to its owning paragraph.\nconst results = body.search(OLD_TAIL, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nconst paragraphs = [];\nfor (let i = 0; i < results.items.length; i++) {\n  paragraphs.push(results.items[i].paragraphs.getFirst());\n}\nawait context.sync();\n\n// Wipe each matched paragraph's runs (text + formatting) and replace with\n// one plain run holding the new sentence, matching the target markup\n// (<w:r><w:t>\u2026</w:t></w:r>, no rPr).\nfor (let i = 0; i < paragraphs.length; i++) {\n  paragraphs[i].clear();\n  paragraphs[i].insertText(NEW_TEXT, Word.InsertLocation.start);\n}\n\nawait context.sync();\n", "ps1": "# The document repeats a \"2018 observation dates for the Perseus\n# constellation\" sentence in four places. Each occurrence is replaced with\n# a single, unformatted run containing the new (translated/updated) Orion\n# sentence -- the \"2018 \" lead-in is dropped and the constellation / date\n# range text is swapped for the Orion one.\n$d = $word.ActiveDocument\n\n$oldTail = \"30 \u039f\u03ba\u03c4\u03c9\u03b2\u03c1\u03af\u03bf\u03c5-8 \u039d\u03bf\u03b5\u03bc\u03b2\u03c1\u03af\u03bf\u03c5 \u03ba\u03b1\u03b9 29 \u039d\u03bf\u03b5\u03bc\u03b2\u03c1\u03af\u03bf\u03c5-8 \u0394\u03b5\u03ba\u03b5\u03bc\u03b2\u03c1\u03af\u03bf\u03c5\"\n$newText = \"\u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2 \u03c0\u03b1\u03c1\u03b1\u03c4\u03ae\u03c1\u03b7\u03c3\u03b7\u03c2 \u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 Orion: 16-25 \u0399\u03b1\u03bd\u03bf\u03c5\u03b1\u03c1\u03af\u03bf\u03c5, 14-23 \u03a6\u03b5\u03b2\u03c1\u03bf\u03c5\u03b1\u03c1\u03af\u03bf\u03c5, 14-24 \u039c\u03b1\u03c1\u03c4\u03af\u03bf\u03c5\"\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = $oldTail\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\n\nwhile ($rng.Find.Execute()) {\n    # Grow the hit to cover its whole owning paragraph (minus the trailing\n    # paragraph mark) so every run in the old sentence gets swept up.\n    $pRng = $rng.Duplicate\n    [void]$pRng.Expand(4)\n    [void]$pRng.MoveEnd(1, -1)\n\n    # Delete the old runs/text, reset formatting so the new run carries no\n    # rPr, then insert the replacement sentence as plain text.\n    [void]$pRng.Delete()\n    [void]$pRng.Font.Reset()\n    [void]$pRng.InsertAfter($newText)\n\n    # Continue searching right after the paragraph we just rewrote.\n    [void]$rng.SetRange($pRng.End, $pRng.End)\n}\n"}
